# Apply the "bang mapping" template data fix:
#  - swap the "DVT - KH" / "So luong - KH" columns (D and E) so that
#    "So luong - KH" now comes before "DVT - KH"
#  - re-style those two columns with an explicit black font colour
#  - move the active selection to K5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap column D and E content (header row 1 + data row 2) -------------
$d1 = $ws.Range("D1").Value()
$e1 = $ws.Range("E1").Value()
$ws.Range("D1").Value = $e1
$ws.Range("E1").Value = $d1

$d2 = $ws.Range("D2").Value()
$e2 = $ws.Range("E2").Value()
$ws.Range("D2").Value = $e2
$ws.Range("E2").Value = $d2

# --- re-style header cells D1:E1 (bold, bordered, explicit black font) ---
$headerRange = $ws.Range("D1:E1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 0
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# --- re-style data cells D2:E2 (explicit black font, no border) ----------
$dataRange = $ws.Range("D2:E2")
$dataRange.Font.Bold = $false
$dataRange.Font.Color = 0

# --- move active selection ------------------------------------------------
$ws.Range("K5").Select()
